$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the power supply description and its unit price.
$ws.Range("B4").Value = "12v 2A Power Supply"
$ws.Range("D4").Value = 5.22

# Move the selection to F5 (matches the author's last selected cell).
$ws.Range("F5").Select()
